# 260127 DataBaseCol: 오타 수정 (insident -> incident) 및 선택 셀 갱신
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# incident 테이블 (행 46~56) - 'insident' 오타를 'incident'로 수정
$ws.Range("B46").Value = "incident"
$ws.Range("C46").Value = "incident_id"
$ws.Range("F46").Value = "incident ID"

$ws.Range("C47").Value = "incident_title"
$ws.Range("F47").Value = "incident 제목"

$ws.Range("C48").Value = "incident_line_name"

$ws.Range("C49").Value = "incident_station_id "

$ws.Range("C50").Value = "incident_station_name"

$ws.Range("C51").Value = "incident_content"

$ws.Range("C52").Value = "incident_status"

# incident_comment 테이블 (행 57~62)
$ws.Range("B57").Value = "incident_comment"
$ws.Range("C58").Value = "incident_id"
$ws.Range("F58").Value = "인시던트 ID ( FK storage.incident(incident_id) )"

# 현재 선택 셀 갱신
$ws.Range("F52").Select() | Out-Null
